$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in row 13 (hour 12) ---
$ws.Range("I13").Value = 19730
$ws.Range("K13").Value = 108706
$ws.Range("R13").Value = 19731.973
$ws.Range("T13").Value = 108716.8706

# --- Update existing values in row 14 (hour 13) ---
$ws.Range("I14").Value = 19881
$ws.Range("K14").Value = 108006
$ws.Range("R14").Value = 19882.9881
$ws.Range("T14").Value = 108016.8006

# --- New row 15 (hour 14) ---
# Column A holds a plain text date (matches the rest of the sheet, which
# stores "2024-08-06" as literal text rather than a date serial). Copy an
# existing text cell down instead of assigning the string directly so
# Excel doesn't auto-convert it into a date value/format.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial()
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 6227
$ws.Range("E15").Value = 19547
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 14601
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 20973
$ws.Range("J15").Value = 22230
$ws.Range("K15").Value = 83578
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 6227.6227
$ws.Range("N15").Value = 19548.9547
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 14602.4601
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 20975.0973
$ws.Range("S15").Value = 22232.223
$ws.Range("T15").Value = 83586.3578

# --- New row 16 (hour 15) ---
$ws.Range("A14").Copy()
$ws.Range("A16").PasteSpecial()
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 6249
$ws.Range("E16").Value = 19436
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 14422
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 20862
$ws.Range("J16").Value = 22243
$ws.Range("K16").Value = 83212
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 6249.6249
$ws.Range("N16").Value = 19437.9436
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 14423.4422
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 20864.0862
$ws.Range("S16").Value = 22245.2243
$ws.Range("T16").Value = 83220.32120000001

# --- New row 17 (hour 16) ---
$ws.Range("A14").Copy()
$ws.Range("A17").PasteSpecial()
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 5989
$ws.Range("E17").Value = 18881
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 14102
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 19536
$ws.Range("J17").Value = 21365
$ws.Range("K17").Value = 79873
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 5989.5989
$ws.Range("N17").Value = 18882.8881
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 14103.4102
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 19537.9536
$ws.Range("S17").Value = 21367.1365
$ws.Range("T17").Value = 79880.98730000001
